$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 39
$ws.Range("H39").Value = 449
$ws.Range("I39").Value = 251.28572
$ws.Range("J39").Value = 504.36
$ws.Range("K39").Value = 753.85716
$ws.Range("L39").Value = 1513.08
$ws.Range("M39").Value = -457.85716
$ws.Range("N39").Value = -2105.08

# row 40
$ws.Range("H40").Value = 3297.0605
$ws.Range("I40").Value = 2563.3157
$ws.Range("K40").Value = 2563.3157
$ws.Range("M40").Value = -2388.3157

# row 42
$ws.Range("H42").Value = 113.1
$ws.Range("I42").Value = 117.75
$ws.Range("K42").Value = 353.25
$ws.Range("M42").Value = -123.25

# row 70
$ws.Range("H70").Value = 1790.6666
$ws.Range("I70").Value = 959.2
$ws.Range("J70").Value = 2384.5715
$ws.Range("K70").Value = 2877.6
$ws.Range("L70").Value = 7153.7145
$ws.Range("M70").Value = -2607.6
$ws.Range("N70").Value = -7693.7145

# row 73
$ws.Range("H73").Value = 1790.6666
$ws.Range("I73").Value = 959.2
$ws.Range("J73").Value = 2384.5715
$ws.Range("K73").Value = 2877.6
$ws.Range("L73").Value = 7153.7145
$ws.Range("M73").Value = -1941.6
$ws.Range("N73").Value = -9025.7145

# row 132
$ws.Range("H132").Value = 1519.2333
$ws.Range("I132").Value = 1556.3572
$ws.Range("K132").Value = 4669.071599999999
$ws.Range("M132").Value = -2139.071599999999

# row 138
$ws.Range("H138").Value = 3001.5264
$ws.Range("I138").Value = 1673.5
$ws.Range("J138").Value = 3355.6667
$ws.Range("K138").Value = 5020.5
$ws.Range("L138").Value = 10067.0001
$ws.Range("M138").Value = 119.5
$ws.Range("N138").Value = -20347.0001

$ws = $wb.Worksheets.Item("ARM")
# row 3
$ws.Range("H3").Value = 7332.1665
$ws.Range("I3").Value = 6164.3335
$ws.Range("J3").Value = 8500
$ws.Range("K3").Value = 6164.3335
$ws.Range("L3").Value = 8500
$ws.Range("M3").Value = -6049.3335
$ws.Range("N3").Value = -8730

# row 31
$ws.Range("H31").Value = 45725.75
$ws.Range("I31").Value = 4249.75
$ws.Range("J31").Value = 87201.75
$ws.Range("K31").Value = 4249.75
$ws.Range("L31").Value = 87201.75
$ws.Range("M31").Value = -3955.75
$ws.Range("N31").Value = -87789.75

$ws = $wb.Worksheets.Item("BSM")
# row 98
$ws.Range("H98").Value = 108000
$ws.Range("J98").Value = 108000
$ws.Range("L98").Value = 108000
$ws.Range("N98").Value = -113990

# row 102
$ws.Range("H102").Value = 75599.2
$ws.Range("I102").Value = 42776.75
$ws.Range("J102").Value = 97480.836
$ws.Range("K102").Value = 42776.75
$ws.Range("L102").Value = 97480.836
$ws.Range("M102").Value = -39531.75
$ws.Range("N102").Value = -103970.836

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# row 31
$ws.Range("H31").Value = 713326.2
$ws.Range("I31").Value = 12210.643
$ws.Range("J31").Value = 1229937.6
$ws.Range("K31").Value = 12210.643
$ws.Range("L31").Value = 1229937.6
$ws.Range("M31").Value = -11915.643
$ws.Range("N31").Value = -1230527.6

# row 34
$ws.Range("H34").Value = 713326.2
$ws.Range("I34").Value = 12210.643
$ws.Range("J34").Value = 1229937.6
$ws.Range("K34").Value = 12210.643
$ws.Range("L34").Value = 1229937.6
$ws.Range("M34").Value = -12008.643
$ws.Range("N34").Value = -1230341.6

# row 92
$ws.Range("H92").Value = 44500
$ws.Range("J92").Value = 44500
$ws.Range("L92").Value = 44500
$ws.Range("N92").Value = -49492

# row 99
$ws.Range("H99").Value = 4005.5
$ws.Range("I99").Value = 3012
$ws.Range("K99").Value = 3012
$ws.Range("M99").Value = -1514

# row 105
$ws.Range("H105").Value = 2172.8333
$ws.Range("I105").Value = 1342.3334
$ws.Range("K105").Value = 1342.3334
$ws.Range("M105").Value = 404.6666

# row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# row 126
$ws.Range("H126").Value = 4005.5
$ws.Range("I126").Value = 3012
$ws.Range("K126").Value = 9036
$ws.Range("M126").Value = -6566

# row 134
$ws.Range("H134").Value = 558430.75
$ws.Range("I134").Value = 910525.4
$ws.Range("J134").Value = 5139.143
$ws.Range("K134").Value = 2731576.2
$ws.Range("L134").Value = 15417.429
$ws.Range("M134").Value = -2729041.2
$ws.Range("N134").Value = -20487.429

$ws = $wb.Worksheets.Item("CUL")
# row 56
$ws.Range("H56").Value = 11423.429
$ws.Range("I56").Value = 11423.429
$ws.Range("K56").Value = 11423.429
$ws.Range("M56").Value = -10893.429

# row 107
$ws.Range("H107").Value = 878.1667
$ws.Range("J107").Value = 878.1667
$ws.Range("L107").Value = 2634.5001
$ws.Range("N107").Value = -6474.5001

# row 113
$ws.Range("H113").Value = 1189.0588
$ws.Range("I113").Value = 618.8570999999999
$ws.Range("J113").Value = 1588.2
$ws.Range("K113").Value = 1856.5713
$ws.Range("L113").Value = 4764.6
$ws.Range("M113").Value = 313.4287000000002
$ws.Range("N113").Value = -9104.6

# row 132
$ws.Range("H132").Value = 1648.72
$ws.Range("I132").Value = 1731.875
$ws.Range("K132").Value = 15586.875
$ws.Range("M132").Value = -13056.875

$ws = $wb.Worksheets.Item("GSM")
# row 3
$ws.Range("H3").Value = 100539.8
$ws.Range("I3").Value = 250125
$ws.Range("J3").Value = 816.3333
$ws.Range("K3").Value = 250125
$ws.Range("L3").Value = 816.3333
$ws.Range("M3").Value = -250009
$ws.Range("N3").Value = -1048.3333

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 4891.6665
$ws.Range("I40").Value = 2349
$ws.Range("J40").Value = 5400.2
$ws.Range("K40").Value = 2349
$ws.Range("L40").Value = 5400.2
$ws.Range("M40").Value = -2213
$ws.Range("N40").Value = -5672.2

# row 50
$ws.Range("H50").Value = 29998
$ws.Range("I50").Value = 29998
$ws.Range("K50").Value = 29998
$ws.Range("M50").Value = -29361

# row 111
$ws.Range("H111").Value = 94744.75
$ws.Range("I111").Value = 87989
$ws.Range("J111").Value = 96996.664
$ws.Range("K111").Value = 87989
$ws.Range("L111").Value = 96996.664
$ws.Range("M111").Value = -83899
$ws.Range("N111").Value = -105176.664

$ws = $wb.Worksheets.Item("WVR")
# row 8
$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2000
$ws.Range("N8").Value = -2280
$ws.Range("M8").ClearContents()

# row 11
$ws.Range("H11").Value = 9999.5
$ws.Range("J11").Value = 9999.5
$ws.Range("L11").Value = 9999.5
$ws.Range("N11").Value = -10283.5

# row 50
$ws.Range("H50").Value = 44999
$ws.Range("J50").Value = 44999
$ws.Range("L50").Value = 44999
$ws.Range("N50").Value = -46261

# row 81
$ws.Range("H81").Value = 500
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 61
$ws.Range("N81").ClearContents()

# row 84
$ws.Range("H84").Value = 500
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 304
$ws.Range("N84").ClearContents()

# row 86
$ws.Range("H86").Value = 57996.668
$ws.Range("J86").Value = 57996.668
$ws.Range("L86").Value = 57996.668
$ws.Range("N86").Value = -60242.668

# row 89
$ws.Range("H89").Value = 57996.668
$ws.Range("J89").Value = 57996.668
$ws.Range("L89").Value = 289983.34
$ws.Range("N89").Value = -301215.34
